$d = $word.ActiveDocument

# 1) Subheading: "Statistical Anomaly Detection for Trafficking & Organized Crime"
#    -> "Geospatial Crime Pattern Analysis | 41,200 NamUs Cases"
$rng = $d.Content
$rng.Find.Execute("Statistical Anomaly Detection for Trafficking & Organized Crime")
if ($rng.Find.Found) {
    $rng.Text = "Geospatial Crime Pattern Analysis | 41,200 NamUs Cases"
}

# 2) First bullet under that subheading: reworded.
$rng2 = $d.Content
$rng2.Find.Execute("Analyzed 41,200 cases across 101 years identifying trafficking corridors at up to 46.86" + [char]0x03C3 + " significance")
if ($rng2.Find.Found) {
    $rng2.Text = "7 statistical methods + 3 ML models detecting trafficking corridors at up to 46.86" + [char]0x03C3 + " significance"
}

# 3) Insert a brand-new bullet paragraph right after the bullet edited above,
#    before the "Built 7-page interactive..." bullet, matching its paragraph
#    formatting (ListParagraph style + same numbering).
$rng3 = $d.Content
$rng3.Find.Execute("7 statistical methods + 3 ML models detecting trafficking corridors at up to 46.86" + [char]0x03C3 + " significance")
$para = $rng3.Paragraphs(1)
$insertPoint = $para.Range
$insertPoint.Collapse(0)
$insertPoint.InsertParagraphAfter()
$newPara = $para.Next()
$newPara.Range.Text = [char]0x2022 + " I-35 corridor acceleration: 170% increase in missing persons, structural break at 2020"

# 4) Last bullet in that block: reworded (contains a straight apostrophe that
#    must NOT be smart-quoted, so assign Range.Text directly instead of using
#    Find.Execute's Replace, which runs AutoCorrect/AutoFormat-as-you-type).
$rng4 = $d.Content
$rng4.Find.Execute("Built 7-page interactive Streamlit dashboard with geospatial visualization")
if ($rng4.Find.Found) {
    $rng4.Text = "Live Streamlit dashboard with spatial autocorrelation (Moran's I), ARIMA forecasting, and LISA clustering"
}

Write-Output "done"
